$d = $word.ActiveDocument

$old = " Escenarios que poseían bajas calificaciones en sus resultados, ahora tienen la calificación más alta."
$new = " Escenarios que poseían bajas calificaciones en sus resultados, ahora tienen la calificación más alta (excepto aquel que aún no haya sido implementado o corregido)."

$range = $d.Content
$found = $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)

if (-not $found) {
    throw "Target sentence not found"
}
